$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.842.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "'2.087.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'235.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'59.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.89%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "'0.0792"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").Value = "'2.395.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "'14.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "'5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "'2.089.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "'37.790.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "'6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").Value = "'71.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "'229.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "'170.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +8.68%  "
$ws.Range("D28").Value = "'9.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'19.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("E32").Value = "  +4.11%  "
$ws.Range("D33").Value = "'0.0633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("D34").Value = "'4.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "'2.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "'3.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.21%  "
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").Value = "'99.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "'1.462.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "'1.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'4.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").Value = "'16.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.76%  "
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").Value = "'47.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.28%  "
